$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the existing text that needs to move around (read BEFORE any writes) ---
$tableDesc = $ws.Range("A55").Value2   # IBGE table description (currently row 55)
$urlText   = $ws.Range("A56").Value2   # URL text, currently hyperlinked (row 56)
$ibgeLabel = $ws.Range("A60").Value2   # "IBGE" label (currently row 60)

# Remove the existing hyperlink object attached to A56 before shuffling values
$ws.Hyperlinks.Delete()

# Row 55: was the table description -> becomes a blank (empty string) "source"
# styled cell, same as it already is, so just clear the text itself.
$ws.Range("A55").Value = "'"

# Row 56: now holds the table description text, styled like the other
# "source" rows (italic, no underline/colour -> drop the old HyperLink look)
$ws.Range("A56").Value = $tableDesc
$ws.Range("A56").Font.Italic = $true
$ws.Range("A56").Font.Underline = $false
$ws.Range("A56").Font.Bold = $false

# Row 57 keeps its existing blank value/style (unchanged)

# Row 58: new row holding the (now plain, non-hyperlinked) URL text, same
# "source" (italic) look as the surrounding rows
$ws.Range("A58").Value = $urlText
$ws.Range("A58").Font.Italic = $true
$ws.Range("A58").Font.Underline = $false
$ws.Range("A58").Font.Bold = $false

# Shift the trailing "IBGE" source block down one row (60/61 -> 61/62)
$ws.Range("A60").Clear()

$ws.Range("A61").Value = $ibgeLabel
$ws.Range("A61").Font.Bold = $true
$ws.Range("A61").Font.Italic = $false
$ws.Range("A61").Font.Underline = $false

$ws.Range("A62").Value = "Observatorio da Lei Geral da Micro e Pequena Empresa available at http://www.leigeral.com.br/portal/main.jsp?lumPageId=FF8081812658D379012665B59AC01CE8"
$ws.Range("A62").Font.Italic = $true
$ws.Range("A62").Font.Underline = $false
$ws.Range("A62").Font.Bold = $false
